# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held values like "11-22-2007-08" (a mangled
# concatenation of the game date and season) for every data row; the
# correct value is the plain ISO game date "2007-11-22".
#
# The column is formatted as General/text, so assigning a literal
# "2007-11-22" string via .Value would make Excel auto-detect it as a
# date and silently convert it to a date serial number. To keep it as
# literal text (matching the original inline-string storage) we prefix
# the assignment with an apostrophe (the standard "force text" trick),
# then reset the cell Style back to "Normal" so no stray number-format /
# quote-prefix styling is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = 58  # column BF

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'2007-11-22"
    $cell.Style = "Normal"
}
